$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Bundás', 2010, 'labrador', 'nagy', 'fiú', 'nem', 'lakos')
    ,('Szuszu', 2016, 'csivava', 'kicsi', 'fiú', 'igen', 'foglalt')
    ,('Karcsi', 2020, 'pitbull', 'közepes', 'fiú', 'igen', 'örökbeadott')
    ,('Don', 2014, 'keverék', 'nagy', 'fiú', 'igen', 'lakos')
    ,('Dunya', 2018, 'németjuhász', 'nagy', 'lány', 'igen', 'lakos')
    ,('Lüszi', 2021, 'kaukázusi juhászkutya', 'nagy', 'lány', 'igen', 'örökbeadott')
    ,('Honey', 2020, 'pitbull', 'közepes', 'fiú', 'igen', 'lakos')
    ,('Amy', 2020, 'rottweiler', 'nagy', 'lány', 'igen', 'lakos')
    ,('Beni', 2012, 'keverék', 'kicsi', 'fiú', 'igen', 'lakos')
    ,('Bianka', 2011, 'keverék', 'közepes', 'lány', 'igen', 'lakos')
    ,('Artúr', 2017, 'tacskó', 'kicsi', 'fiú', 'igen', 'lakos')
    ,('Gino', 2011, 'juhászkutya', 'nagy', 'lány', 'igen', 'lakos')
    ,('Lángos', 2020, 'pitbull', 'közepes', 'fiú', 'igen', 'lakos')
    ,('Cserkés', 2020, 'keverék', 'nagy', 'fiú', 'igen', 'lakos')
    ,('Dorisz', 2015, 'pitbull', 'nagy', 'lány', 'igen', 'lakos')
    ,('Cheester', 2018, 'staffordshire', 'közepes', 'fiú', 'igen', 'lakos')
    ,('Artemisz', 2020, 'kaukázusi juhászkutya', 'nagy', 'lány', 'igen', 'örökbeadott')
    ,('Tyson', 2022, 'husky', 'nagy', 'fiú', 'nem', 'foglalt')
    ,('Nudli', 2017, 'tacskó', 'közepes', 'fiú', 'igen', 'foglalt')
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("M9").Select()
